# Adds "Turno", "Rango de duración" and "Intensidad (según BPM)" columns
# (F, G, H) to the caminatas dataset, derived from the existing walk
# records (time of day -> Turno, duration -> Rango de duración,
# average heart rate -> Intensidad).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- column F: Turno ------------------------------------------------------
$ws.Cells.Item(1, 6).Value = "Turno"
$turno = @{
    2  = "Mediodía"; 3  = "Noche"
    4  = "Mediodía"; 5  = "Noche"
    6  = "Mediodía"; 7  = "Noche"
    8  = "Mediodía"; 9  = "Noche"
    10 = "Mediodía"; 11 = "Noche"
    12 = "Mediodía"; 13 = "Noche"
    14 = "Mediodía"; 15 = "Noche"
    16 = "Mediodía"; 17 = "Noche"
    18 = "Mediodía"; 19 = "Noche"
    20 = "Mediodía"; 21 = "Noche"
    22 = "Mediodía"; 23 = "Noche"
    24 = "Mediodía"; 25 = "Noche"
    26 = "Mediodía"; 27 = "Noche"
    28 = "Mediodía"; 29 = "Noche"
    30 = "Mediodía"
}
for ($r = 2; $r -le 30; $r++) {
    $ws.Cells.Item($r, 6).Value = $turno[$r]
}

# --- column G: Rango de duración ------------------------------------------
$ws.Cells.Item(1, 7).Value = "Rango de duración"
$rango = @{
    2  = "Corta"; 3  = "Media"
    4  = "Corta"; 5  = "Media"
    6  = "Media"; 7  = "Larga"
    8  = "Corta"; 9  = "Media"
    10 = "Corta"; 11 = "Larga"
    12 = "Corta"; 13 = "Media"
    14 = "Corta"; 15 = "Larga"
    16 = "Corta"; 17 = "Larga"
    18 = "Corta"; 19 = "Larga"
    20 = "Corta"; 21 = "Larga"
    22 = "Corta"; 23 = "Media"
    24 = "Corta"; 25 = "Media"
    26 = "Corta"; 27 = "Larga"
    28 = "Corta"; 29 = "Larga"
    30 = "Corta"
}
for ($r = 2; $r -le 30; $r++) {
    $ws.Cells.Item($r, 7).Value = $rango[$r]
}

# --- column H: Intensidad (según BPM) -------------------------------------
$ws.Cells.Item(1, 8).Value = "Intensidad (según BPM)"
$intensidad = @{
    2  = "Media"; 3  = "Media"
    4  = "Baja";  5  = "Baja"
    6  = "Baja";  7  = "Media"
    8  = "Baja";  9  = "Baja"
    10 = "Media"; 11 = "Alta"
    12 = "Media"; 13 = "Media"
    14 = "Media"; 15 = "Alta"
    16 = "Media"; 17 = "Alta"
    18 = "Media"; 19 = "Alta"
    20 = "Media"; 21 = "Alta"
    22 = "Baja";  23 = "Alta"
    24 = "Baja";  25 = "Media"
    26 = "Media"; 27 = "Alta"
    28 = "Media"; 29 = "Alta"
    30 = "Media"
}
for ($r = 2; $r -le 30; $r++) {
    $ws.Cells.Item($r, 8).Value = $intensidad[$r]
}

# --- column widths ----------------------------------------------------------
$ws.Columns.Item(4).AutoFit()
$ws.Columns.Item(5).AutoFit()
$ws.Columns.Item(7).AutoFit()
$ws.Columns.Item(8).AutoFit()

# --- selection --------------------------------------------------------------
$ws.Range("H27").Select()
